$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for data rows 2-19
# from 45170 to 45174, preserving existing cell formatting/style.
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
